$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/4/2023  Through  9/10/2023"

# --- Row 15: D15/E15 become shared-string "0" / "***.*" cells, copying style+type from analogous cells ---
$ws.Range("C15").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))

# --- Numeric cell value updates ---
$ws.Range("L14").Value = -40
$ws.Range("N14").Value = -72.727272727272
$ws.Range("M15").Value = -28.571428571428
$ws.Range("N15").Value = -47.368421052631
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 141
$ws.Range("J16").Value = 138
$ws.Range("K16").Value = 2.173913043478
$ws.Range("L16").Value = 17.5
$ws.Range("M16").Value = -8.441558441558
$ws.Range("N16").Value = -56.346749226006
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 41.666666666666
$ws.Range("I17").Value = 210
$ws.Range("J17").Value = 180
$ws.Range("K17").Value = 16.666666666666
$ws.Range("L17").Value = 38.157894736842
$ws.Range("M17").Value = 56.716417910447
$ws.Range("N17").Value = 21.387283236994
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -12.5
$ws.Range("I18").Value = 95
$ws.Range("J18").Value = 73
$ws.Range("K18").Value = 30.136986301369
$ws.Range("L18").Value = 30.136986301369
$ws.Range("M18").Value = -57.777777777777
$ws.Range("N18").Value = -84.477124183006
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -41.666666666666
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 13.636363636363
$ws.Range("I19").Value = 457
$ws.Range("J19").Value = 375
$ws.Range("K19").Value = 21.866666666666
$ws.Range("L19").Value = 44.164037854889
$ws.Range("M19").Value = 38.484848484848
$ws.Range("N19").Value = 51.827242524916
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -12.5
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 32
$ws.Range("H20").Value = -15.625
$ws.Range("I20").Value = 372
$ws.Range("J20").Value = 205
$ws.Range("K20").Value = 81.463414634146
$ws.Range("L20").Value = 85.074626865671
$ws.Range("M20").Value = 149.664429530201
$ws.Range("N20").Value = -74.503084304318
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -21.428571428571
$ws.Range("F21").Value = 120
$ws.Range("G21").Value = 113
$ws.Range("H21").Value = 6.194690265486
$ws.Range("I21").Value = 1288
$ws.Range("J21").Value = 986
$ws.Range("K21").Value = 30.628803245436
$ws.Range("L21").Value = 46.864310148232
$ws.Range("M21").Value = 28.031809145129
$ws.Range("N21").Value = -55.555555555555
$ws.Range("M22").Value = -10
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 200
$ws.Range("L23").Value = 37.5
$ws.Range("M23").Value = 33.333333333333
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 52.941176470588
$ws.Range("F24").Value = 120
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = 18.811881188118
$ws.Range("I24").Value = 988
$ws.Range("J24").Value = 811
$ws.Range("K24").Value = 21.824907521578
$ws.Range("L24").Value = 55.102040816326
$ws.Range("M24").Value = -5.363984674329
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 26.666666666666
$ws.Range("I25").Value = 356
$ws.Range("J25").Value = 314
$ws.Range("K25").Value = 13.375796178343
$ws.Range("L25").Value = 18.666666666666
$ws.Range("M25").Value = 24.912280701754
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 23
$ws.Range("J26").Value = 25
$ws.Range("K26").Value = -8
$ws.Range("L26").Value = 76.923076923076
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 29
$ws.Range("J27").Value = 37
$ws.Range("K27").Value = -21.621621621621
$ws.Range("L27").Value = -3.333333333333
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = -40
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = -53.333333333333
